$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the first 4 data rows (rows 2-5), shifting the rest up.
$ws.Range("A2:B5").EntireRow.Delete()
